$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.718.12"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "2.220.21"
$ws.Range("E3").Value = "  -5.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "`'294.73"
$ws.Range("E5").Value = "  -5.16%  "
$ws.Range("D6").Value = "`'83.91"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "`'0.466"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").Value = "`'0.0784"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").Value = "`'29.70"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "`'47.69"
$ws.Range("E12").Value = "  -9.02%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").Value = "2.564.40"
$ws.Range("E14").Value = "  -5.40%  "
$ws.Range("D15").Value = "`'6.29"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "`'14.12"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "2.221.23"
$ws.Range("E17").Value = "  -6.38%  "
$ws.Range("D18").Value = "`'0.719"
$ws.Range("E18").Value = "  -5.45%  "
$ws.Range("D19").Value = "39.625.69"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "0.0₃0882"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "`'5.75"
$ws.Range("E21").Value = "  -5.96%  "
$ws.Range("D22").Value = "`'65.02"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "`'10.45"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "`'232.22"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "`'22.82"
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").Value = "`'9.18"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "`'32.21"
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("D32").Value = "`'149.74"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").Value = "`'1.00"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "`'4.82"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "`'2.37"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "`'0.0704"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "`'15.99"
$ws.Range("E37").Value = "  +2.37%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "`'0.0971"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "`'2.66"
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("D41").Value = "`'1.65"
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("D42").Value = "`'3.67"
$ws.Range("E42").Value = "  -5.37%  "
$ws.Range("D43").Value = "1.939.73"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("D45").Value = "`'0.0266"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "`'9.44"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "`'16.08"
$ws.Range("E47").Value = "  -9.05%  "
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").Value = "2.434.26"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").Value = "`'70.84"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "`'88.92"
